# Rename the original (only) worksheet from "Hoja1" to "version_1",
# then duplicate it right after itself to create "version_2" — an
# identical copy except the "Dataset" cell in row 2 (E2) which should
# read "smallest (N=100)" instead of "smallest (N=28)".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "version_1"

# Copy version_1 and place the copy immediately after it.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "version_2"

# The only real data change: row 2's dataset size note.
$ws2.Range("E2").Value = "smallest (N=100)"

# The copied sheet also carries forward the wrap-text formatting onto the
# otherwise-empty cells at the end of the summary row, same as the other
# cells in that row/column band.
$ws2.Range("A8").WrapText = $true
$ws2.Range("B8").WrapText = $true
$ws2.Range("G8").WrapText = $true

# Match the recorded selections on each tab.
$ws1.Range("A4:XFD4").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("G4").Select() | Out-Null
